# Insert a new data row at row 59 (Hortaliza / Berenjena weekly price update).
# This shifts the existing rows 59-123 down to 60-124 and extends the sheet
# dimension from A1:R123 to A1:R124, matching the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(59).Insert()

$ws.Range("A59").Value = 6
$ws.Range("B59").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C59").Value = "Metropolitana"
$ws.Range("D59").Value = 44512
$ws.Range("E59").Value = 13
$ws.Range("F59").Value = 100112001
$ws.Range("G59").Value = "Berenjena"
$ws.Range("H59").Value = "Sin especificar"
$ws.Range("I59").Value = "Primera"
$ws.Range("J59").Value = 230
$ws.Range("K59").Value = 12000
$ws.Range("L59").Value = 13000
$ws.Range("M59").Value = 12435
$ws.Range("N59").Value = "`$/caja 60 unidades"
$ws.Range("O59").Value = "Provincia de Huasco"
$ws.Range("P59").Value = 207
$ws.Range("Q59").Value = 60
$ws.Range("R59").Value = "Hortaliza"
